$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.0498220640569395
$ws1.Range("C2").Value = 0.0498220640569395
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 0.09491525423728814
$ws1.Range("F2").Value = 0.2077151335311573
$ws1.Range("G2").Value = 0.5768621236133122
$ws1.Range("H2").Value = 0.7501003210272873
$ws1.Range("I2").Value = 28
$ws1.Range("J2").Value = 534
$ws1.Range("K2").Value = 0
$ws1.Range("L2").Value = 0

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# Row 2 ("0")
$ws2.Range("B2").Value = 0
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 0

# Row 3 ("1")
$ws2.Range("B3").Value = 0.0498220640569395
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.09491525423728814

# Row 4 ("accuracy")
$ws2.Range("B4").Value = 0.0498220640569395
$ws2.Range("C4").Value = 0.0498220640569395
$ws2.Range("D4").Value = 0.0498220640569395
$ws2.Range("E4").Value = 0.0498220640569395

# Row 5 ("macro avg")
$ws2.Range("B5").Value = 0.02491103202846975
$ws2.Range("C5").Value = 0.5
$ws2.Range("D5").Value = 0.04745762711864407

# Row 6 ("weighted avg")
$ws2.Range("B6").Value = 0.002482238066893783
$ws2.Range("C6").Value = 0.0498220640569395
$ws2.Range("D6").Value = 0.004728873876590867

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 ("Actual 0")
$ws3.Range("B2").Value = 0
$ws3.Range("C2").Value = 534

# Row 3 ("Actual 1")
$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 28
